# Auto-generated edit script for cryptos.xlsx update
# Commit: Updated cryptos list on Fri Jul  7 21:56:19 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) values: force text format to preserve exact formatting ---
# (avoids Excel auto-converting numeric-looking strings to floats, which would
#  drop significant trailing zeros like "165.40" -> 165.4)
$dCells = @("D2", "D3", "D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D19", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D33", "D36", "D39", "D40", "D42", "D43", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($ref in $dCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# --- Apply new cell values ---
$ws.Range("D2").Value = "30.331.07"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "1.870.82"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "235.57"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "0.4676"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "0.2845"
$ws.Range("E8").Value = "  +0.65%  "
$ws.Range("D9").Value = "0.06533"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").Value = "21.39"
$ws.Range("E10").Value = "  +4.54%  "
$ws.Range("D11").Value = "0.07871"
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("D12").Value = "97.96"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "1.866.09"
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("D14").Value = "5.098"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").Value = "0.6757"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").Value = "278.97"
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("D17").Value = "30.315.59"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").Value = "5.514"
$ws.Range("E19").Value = "  +2.40%  "
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.000007310"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.106.77"
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "6.158"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").Value = "165.40"
$ws.Range("E25").Value = "  -1.35%  "
$ws.Range("E26").Value = "  -1.84%  "
$ws.Range("D27").Value = "19.11"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("D28").Value = "1.933"
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("D29").Value = "1.378"
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("D30").Value = "0.09630"
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").Value = "4.375"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("D33").Value = "4.090"
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("E35").Value = "  +3.40%  "
$ws.Range("D36").Value = "0.7066"
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("D39").Value = "6.281"
$ws.Range("E39").Value = "  -3.54%  "
$ws.Range("D40").Value = "2.528"
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("E41").Value = "  +2.85%  "
$ws.Range("D42").Value = "1.944"
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("D43").Value = "0.8501"
$ws.Range("E43").Value = "  -1.34%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D46").Value = "103.92"
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("D47").Value = "7.168"
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("D48").Value = "9.196"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("D49").Value = "935.43"
$ws.Range("E49").Value = "  -4.93%  "
$ws.Range("D50").Value = "34.15"
$ws.Range("E50").Value = "  +0.98%  "
$ws.Range("D51").Value = "0.1125"
$ws.Range("E51").Value = "  -1.80%  "

Write-Output "Applied 90 cell updates"
